# Legend slide (slide 1) touch-ups:
#  1. Thicken the grey "tick absent" connector line (12pt -> 20pt).
#  2. Thicken the orange connector line (12pt -> 20pt) and make it solid
#     instead of dashed.
#  3. Resize/reposition the orange triangle marker and lock its aspect
#     ratio.
#  4. Add a new grey oval marker (duplicated from the triangle, then
#     converted to an oval) positioned near the top of the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Grey connector that lives inside the "Group 3" group ---------------
$group = $s.Shapes.Item(1)
$greyConnector = $group.GroupItems.Item(1)
$greyConnector.Line.Weight = 20   # 20pt * 12700 EMU/pt = 254000 EMU

# --- 2) Orange dashed connector --------------------------------------------
$orangeConnector = $s.Shapes.Item(2)
$orangeConnector.Line.Weight = 20      # 254000 EMU
$orangeConnector.Line.DashStyle = 1    # msoLineSolid

# --- 3) Orange triangle marker ---------------------------------------------
$triangle = $s.Shapes.Item(4)
$triangle.LockAspectRatio = -1   # msoTrue -> adds <a:spLocks noChangeAspect="1"/>
# New position/size (EMU 6499192,4463176 / 1116000x1115053) expressed in
# points (1 pt = 12700 EMU) with enough precision to survive the host's
# internal float32 storage round-trip.
$triangle.Left = 511.7474015748032
$triangle.Top = 351.4311811023622
$triangle.Width = 87.8740157480315
$triangle.Height = 87.79945081889764

# --- 4) New grey oval marker -------------------------------------------
# Duplicate the (already restyled) triangle so the new shape inherits the
# same p:style/p:txBody boilerplate, then turn it into an oval and
# reposition/recolor it.
$oval = $triangle.Duplicate()
$oval.Name = "Oval 1"
$oval.AutoShapeType = 9   # msoShapeOval
$oval.LockAspectRatio = -1

# New position/size (EMU 6499192,1051558 / 1116000x1116000) in points.
$oval.Left = 511.7474015748032
$oval.Top = 82.79984251968504
$oval.Width = 87.8740157480315
$oval.Height = 87.8740157480315

$oval.Fill.ForeColor.RGB = 13027014   # 0xC6C6C6
